# Update "Datos actualizados a ..." timestamp in the title cell (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 04:35"

# Honduras overtakes Sudan in the ranking (sheet sorted descending by Casos
# totales / column B), so the two rows swap country + Sudan keeps its old
# figures while Honduras receives its updated figures.
$ws.Range("A73").Value = "Honduras"
$ws.Range("B73").Value = 3204
$ws.Range("C73").Value = 104
$ws.Range("D73").Value = 397
$ws.Range("E73").Value = 2651
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 156

$ws.Range("A74").Value = "Sudan"
$ws.Range("B74").Value = 3138
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 309
$ws.Range("E74").Value = 2708
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 121

# Sri Lanka (row 105) updated figures
$ws.Range("B105").Value = 1055
$ws.Range("C105").Value = 7
$ws.Range("E105").Value = 442

# Trinidad y Tobago (row 168) updated figures
$ws.Range("D168").Value = 108
$ws.Range("E168").Value = 0
